$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.289.73'
$ws.Range('E2').Value = '  -3.03%  '
$ws.Range('D3').Value = '2.224.42'
$ws.Range('E3').Value = '  -6.03%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '296.94'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '82.87'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.511'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.92%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.469'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0777'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '29.22'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.83'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -8.89%  '
$ws.Range('E13').Value = '  -1.85%  '
$ws.Range('D14').Value = '2.556.12'
$ws.Range('E14').Value = '  -6.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.23'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.11'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.89%  '
$ws.Range('D17').Value = '2.222.55'
$ws.Range('E17').Value = '  -4.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.716'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.13%  '
$ws.Range('D19').Value = '39.120.48'
$ws.Range('E19').Value = '  -3.21%  '
$ws.Range('D20').Value = '0.0₃0873'
$ws.Range('E20').Value = '  -3.70%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.73'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '64.78'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.20'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '227.30'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.82%  '
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.40'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.77'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.61'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.17'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.10'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.53%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '31.97'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.78%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '147.97'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.86'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0699'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.32'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.110'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.68'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0959'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '14.83'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.63'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.79%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.71'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.53%  '
$ws.Range('D43').Value = '1.913.32'
$ws.Range('E43').Value = '  -2.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0258'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.79%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.01'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -16.00%  '
$ws.Range('E46').Value = '  -4.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.64'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '16.01'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -9.36%  '
$ws.Range('D49').Value = '2.426.39'
$ws.Range('E49').Value = '  -6.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '70.56'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '87.23'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.99%  '
